# Adds match-statistics columns (S, SC, QS, QSC, BTTS%, Blank%, CS%, FG%, FC%)
# to both the "home" and "away" sheets of the Ligue1 Home/Away workbook, for
# all 19 teams (rows 2-19). Mirrors the source dataframe pipeline output.
$wb = $excel.ActiveWorkbook

$cols = @("J", "K", "L", "M", "N", "O", "P", "Q", "R")
$headers = @("S", "SC", "QS", "QSC", "BTTS%", "Blank%", "CS%", "FG%", "FC%")

# home: one row per team (rows 2-19), columns J-R as listed in $cols
$data_home = @(
    @(276, 120, 62, 21, 0.47, 0, 0.53, 0.68, 0.32),
    @(278, 203, 58, 36, 0.58, 0.11, 0.32, 0.63, 0.37),
    @(322, 207, 47, 36, 0.67, 0.14, 0.24, 0.62, 0.33),
    @(297, 215, 41, 24, 0.5, 0.1, 0.5, 0.5, 0.4),
    @(312, 215, 56, 31, 0.57, 0.14, 0.29, 0.67, 0.33),
    @(310, 212, 67, 24, 0.65, 0.15, 0.3, 0.45, 0.45),
    @(316, 200, 44, 28, 0.38, 0.29, 0.57, 0.52, 0.24),
    @(318, 217, 47, 29, 0.48, 0.14, 0.43, 0.67, 0.29),
    @(286, 248, 55, 38, 0.6, 0.15, 0.35, 0.35, 0.55),
    @(257, 223, 44, 30, 0.6, 0.2, 0.25, 0.5, 0.45),
    @(277, 210, 38, 32, 0.47, 0.32, 0.26, 0.42, 0.53),
    @(320, 290, 49, 31, 0.55, 0.2, 0.45, 0.4, 0.4),
    @(222, 310, 33, 30, 0.52, 0.33, 0.29, 0.43, 0.43),
    @(148, 127, 19, 22, 0.45, 0.45, 0.27, 0.36, 0.45),
    @(244, 247, 35, 38, 0.4, 0.45, 0.25, 0.25, 0.65),
    @(275, 263, 32, 51, 0.38, 0.43, 0.24, 0.33, 0.62),
    @(262, 250, 39, 24, 0.65, 0.3, 0.15, 0.35, 0.55),
    @(127, 110, 9, 15, 0.45, 0.55, 0.09, 0.27, 0.64)
)

# away: one row per team (rows 2-19), columns J-R as listed in $cols
$data_away = @(
    @(274, 298, 59, 44, 0.63, 0.05, 0.37, 0.89, 0.05),
    @(282, 282, 49, 45, 0.71, 0.19, 0.14, 0.52, 0.43),
    @(275, 256, 46, 42, 0.43, 0.24, 0.48, 0.48, 0.38),
    @(276, 229, 39, 25, 0.37, 0.37, 0.37, 0.53, 0.37),
    @(234, 244, 23, 40, 0.5, 0.3, 0.3, 0.55, 0.35),
    @(212, 258, 42, 44, 0.65, 0.2, 0.2, 0.45, 0.5),
    @(237, 329, 37, 47, 0.48, 0.38, 0.29, 0.52, 0.33),
    @(245, 246, 23, 33, 0.45, 0.35, 0.4, 0.4, 0.4),
    @(233, 276, 32, 41, 0.4, 0.3, 0.3, 0.55, 0.45),
    @(244, 275, 36, 42, 0.5, 0.35, 0.25, 0.5, 0.4),
    @(267, 235, 36, 36, 0.55, 0.35, 0.15, 0.4, 0.55),
    @(210, 321, 27, 45, 0.6, 0.3, 0.2, 0.4, 0.5),
    @(275, 223, 40, 44, 0.57, 0.29, 0.33, 0.48, 0.33),
    @(191, 237, 22, 45, 0.53, 0.32, 0.26, 0.37, 0.53),
    @(94, 227, 13, 31, 0.5, 0.33, 0.17, 0.42, 0.58),
    @(221, 233, 29, 47, 0.62, 0.24, 0.14, 0.33, 0.67),
    @(109, 201, 10, 30, 0.5, 0.5, 0.25, 0.17, 0.58),
    @(173, 323, 26, 41, 0.55, 0.45, 0.1, 0.4, 0.5)
)

$sheetData = @{ "home" = $data_home; "away" = $data_away }

foreach ($sheetName in @("home", "away")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Header row: write labels, then copy the existing bold/border/center style from A1
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + "1").Value = $headers[$c]
    }
    $ws.Range("A1").Copy()
    $ws.Range("J1:R1").PasteSpecial(-4122)  # xlPasteFormats

    # Data rows 2-19
    $rows = $sheetData[$sheetName]
    for ($r = 0; $r -lt $rows.Length; $r++) {
        $rowValues = $rows[$r]
        $rowNum = $r + 2
        for ($c = 0; $c -lt $cols.Length; $c++) {
            $ws.Range($cols[$c] + $rowNum).Value = $rowValues[$c]
        }
    }
}

Write-Host "Added S/SC/QS/QSC/BTTS%/Blank%/CS%/FG%/FC% columns to home and away sheets"